$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns I ("I0") and J ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing H1 header cell (bold,
# bordered, centered/top-aligned) by copying its format onto I1:J1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows 2-19 for columns I and J ---
$data = @(
    @(4, 9),
    @(6, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 7),
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
